# Apply the "Top 10 communities" restructuring described by the commit.
#
# The sheet grows from 10 data rows (A2:A11) to 12 data rows (A2:A13):
#   - the combined "Chromatin remodelling .../Transcriptional regulation ..."
#     entry is split into two separate rows
#   - the combined "RNA metabolism (... / ...)" entry is split into two
#     separate rows
#   - "Endothelial barrier function and viral entry modulation" is replaced
#     with a re-capitalised variant "Endothelial Barrier Function and Viral
#     Entry Modulation"
#
# New text values must be written in a specific first-use order so that the
# shared-string table receives them in the same order the canonical file
# uses (sharedStrings entries are appended in first-write order).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top 10 communities")

# --- Phase 1 -----------------------------------------------------------
# Write every brand-new string value for the first time, in the exact
# order they must appear in the shared string table.
$ws.Range("A5").Value2  = "Endothelial Barrier Function and Viral Entry Modulation"
$ws.Range("A3").Value2  = "Transcriptional regulation and chromatin remodelling"
$ws.Range("A2").Value2  = "Chromatin remodelling and transcriptional regulation"
$ws.Range("A11").Value2 = "RNA metabolism and modification in viral infection"
$ws.Range("A12").Value2 = "RNA metabolism and viral defense mechanism"

# --- Phase 2 -----------------------------------------------------------
# Fill in the remaining rows, which reuse already-existing shared strings.
$ws.Range("A4").Value2  = "DNA repair and maintenance in response to viral infection"
$ws.Range("A6").Value2  = "Extracellular matrix organization and cell adhesion in response to viral infection"
$ws.Range("A7").Value2  = "Kynurenine pathway modulation and immune response regulation"
$ws.Range("A8").Value2  = "Oxidative stress response and protein quality control in neurodegeneration"
$ws.Range("A9").Value2  = "Protein quality control and intracellular signalling"
$ws.Range("A10").Value2 = "RAS pathway modulation and apoptosis regulation"
$ws.Range("A13").Value2 = "Urea cycle and redox homeostasis"

# Give the new row (A13) the same look as the rest of the data rows
# (left-aligned body style, matching A2:A12).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value2 = "Urea cycle and redox homeostasis"
$excel.CutCopyMode = 0

# --- AutoFilter ----------------------------------------------------------
# Grow the autofilter range to cover the two new rows.
$ws.AutoFilterMode = $false
$ws.Range("A1:A12").AutoFilter() | Out-Null

# --- Defined name ----------------------------------------------------------
# The hidden _FilterDatabase name for this sheet must track the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*" -and $n.RefersTo -like "*Top 10 communities*") {
        $n.RefersTo = "='Top 10 communities'!`$A`$1:`$A`$12"
    }
}

# --- Selection (cosmetic, matches the saved view state) --------------------
$ws.Range("A9").Select() | Out-Null

Write-Host "edit complete"
